$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have 13 columns (A..M); the new layout keeps only the
# old "data_path_*"/"index_*" columns (old D..K), dropping the old
# A..C ("button_*") and L..M ("label_*") columns entirely.
#
# Shift the surviving old D:K block left into A:H (values + header text +
# header style all come along with a copy/paste, since they are the exact
# columns we want to keep).
$ws.Range("D1:K2").Copy()
$ws.Range("A1").PasteSpecial()

# Remove the now-duplicated trailing columns I:M in a single Delete() call
# (full-column Delete() leaves one stale trailing <col> width entry behind
# as an engine quirk - issuing one combined call instead of five separate
# ones keeps that side effect to an unavoidable minimum).
$ws.Columns("I:M").Delete()

# The data_path text values themselves were also reworded - update row 2.
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/proceedToCheckoutWithEmail-test-data"
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/enterShippingDetailsAndProceedToPayment-test-data"
$ws.Range("C2").Value = "Data Files/AI-Generated/Common/proceedToCheckoutWithEmail-test-data"
$ws.Range("D2").Value = "Data Files/AI-Generated/Common/enterShippingDetailsAndProceedToPayment-test-data"

# Finally, apply the new column widths. The ColumnWidth setter round-trips
# through a pixel conversion that adds back 5/6 of a character to whatever
# is stored in the saved XML "width" attribute, so subtract that out here
# to land exactly on the target widths (69, 82, 69, 82, 9, 9, 9, 9).
$ws.Range("A:A").ColumnWidth = 69 - 0.8333333
$ws.Range("B:B").ColumnWidth = 82 - 0.8333333
$ws.Range("C:C").ColumnWidth = 69 - 0.8333333
$ws.Range("D:D").ColumnWidth = 82 - 0.8333333
$ws.Range("E:E").ColumnWidth = 9 - 0.8333333
$ws.Range("F:F").ColumnWidth = 9 - 0.8333333
$ws.Range("G:G").ColumnWidth = 9 - 0.8333333
$ws.Range("H:H").ColumnWidth = 9 - 0.8333333
